# Apply the edits described by the commit "modified test data file":
#  - A1 shared string changes from "Sup" to "Supriya"
#  - A new hyperlinked cell B1 is added with display text/value "Sup@123"
#    (styled with the built-in "Hyperlink" cell style, matching C1's
#    pre-existing style)
#  - The worksheet selection is narrowed from B1:E2 down to just B1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the existing text in A1.
$ws.Range("A1").Value = "Supriya"

# 2. Add B1 with its hyperlink (creates the shared string "Sup@123",
#    the <hyperlinks> entry and the external relationship).
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:Sup@123", "", "", "Sup@123")

# Make sure B1 carries the same "Hyperlink" cell style already used by C1.
$ws.Range("B1").Style = "Hyperlink"

# 3. Collapse the selection so only B1 is selected/active.
$null = $ws.Range("B1").Select()
